# Update the last data row (row 32, year 2020) with refreshed figures
# ("Actualización desde MV -datos-")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B32").Value = 44263751
$ws.Range("C32").Value = 35697172
$ws.Range("F32").Value = 357809
$ws.Range("G32").Value = 891178
$ws.Range("H32").Value = 838970
$ws.Range("I32").Value = 2354488
$ws.Range("J32").Value = 51513332
$ws.Range("K32").Value = 15081661
$ws.Range("L32").Value = 6179012
$ws.Range("M32").Value = 1938521
$ws.Range("N32").Value = 19242319
$ws.Range("O32").Value = 8900542
$ws.Range("P32").Value = 171277
$ws.Range("Q32").Value = -7249582
$ws.Range("R32").Value = 7025832
$ws.Range("T32").Value = 4041888
$ws.Range("U32").Value = 3011804
$ws.Range("V32").Value = 44291610
$ws.Range("W32").Value = 58567024
$ws.Range("X32").Value = -14275414
